# "se estima el error fault al publicar"
# - Replace the placeholder test project names with the real published
#   project id, clear the now-unused second row, and stamp the same
#   "error/fault" highlight style (white-on-black Segoe UI, thin grey
#   border, wrapped text) down through a couple of extra blank rows so the
#   sheet has room to grow when publishing runs again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content -----------------------------------------------------
# A2 held "TestRuben" -> now holds the real promo/project id.
$ws.Range("A2").Value = "AmxCoFijPosPromo023_B26"

# A3 held "TestRuben1" -> no longer needed, blank it out (keeps the cell,
# drops the shared-string entry).
$ws.Range("A3").ClearContents() | Out-Null

# --- Formatting ---------------------------------------------------------
# Apply the "fault" highlight formatting to A2:A5 (this also materialises
# A4/A5 as real, styled, empty cells, extending the sheet dimension).
$fmt = $ws.Range("A2:A5")

$fmt.Font.Name = "Segoe UI"
$fmt.Font.Size = 10
$fmt.Font.Color = 0

$fmt.Interior.Color = 16777215
$fmt.Interior.PatternColor = 0

$fmt.Borders.Color = 13750737

$fmt.WrapText = $true

# --- View state -----------------------------------------------------
$ws.Range("D13").Select() | Out-Null
